$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.509.89"
$ws.Range("E2").Value = "  -2.78%  "
$ws.Range("D3").Value = "2.414.91"
$ws.Range("E3").Value = "  +6.72%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "2.782.54"
$ws.Range("E14").Value = "  +6.88%  "
$ws.Range("D15").Value = "2.418.07"
$ws.Range("E15").Value = "  +6.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.838"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.95%  "
$ws.Range("D18").Value = "45.373.90"
$ws.Range("E18").Value = "  -3.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("D20").Value = "0.0₃0941"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.22%  "
$ws.Range("E31").Value = "  +14.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "148.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0766"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +17.13%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.14%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").Value = "1.995.00"
$ws.Range("E42").Value = "  +12.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "88.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +25.80%  "
$ws.Range("E48").Value = "  +8.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "100.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.66%  "
$ws.Range("D50").Value = "2.651.37"
$ws.Range("E50").Value = "  +6.86%  "
$ws.Range("E51").Value = "  -1.23%  "
